# Strip the erroneous trailing "16" from the verse references in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $value = $cell.Value2
    if ($value -ne $null -and $value.EndsWith("16")) {
        $cell.Value2 = $value.Substring(0, $value.Length - 2)
    }
}
